# Adds new residential-area nodes (CVA, RA, PA, SA, EA) to the distance dataset:
#   - Cougar Village Apartments (CVA)
#   - Reserves Student Apartments (RA)
#   - Parc at 720 Student Apartments (PA)
#   - The Social Student Apartments (SA)
#   - Enclave Student Apartments (EA)
# Extends the "Miles" distance table (cols R:V, rows 11-23) with the new nodes,
# relabels the explanatory notes, and adjusts the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Move/Update NOTE block (old B19/C19:C22 -> new A18/B18:B20) ----
$ws.Range('C19').ClearContents()
$ws.Range('C20').ClearContents()
$ws.Range('C21').ClearContents()
$ws.Range('C22').ClearContents()
$ws.Range('B19').ClearContents()

$ws.Range('A18').Value = 'NOTE:'
$ws.Range('B18').Value = 'Miles Distance Table Begins on G10 (with residential areas)'
$ws.Range('B19').Value = 'Minutes Distance Table Begins on BA9 (no residential areas)'
$ws.Range('B20').Value = 'Location Annotations and Node Coordinates Begins on A1'

# ---- New location annotation + coordinate labels (I3:I7) ----
$ws.Range('I3').Value = 'Cougar Village Apartments = CVA (38.80438098482753, -89.99518772883548)'
$ws.Range('I4').Value = 'Reserves Student Apartments = RA (38.80518785768586, -90.00950273315568)'
$ws.Range('I5').Value = 'Parc at 720 Student Apartments = PA (38.80370049072766, -90.0118599711942)'
$ws.Range('I6').Value = 'The Social Student Apartments = SA (38.78475200026057, -89.98201421937024)'
$ws.Range('I7').Value = 'Enclave Student Apartments = EA (38.80546832836268, -89.97140624030128)'

# ---- Minutes table note (BB9) ----
$ws.Range('BB9').Value = '(NO RESIDENTIAL AREAS)'

# ---- Miles table new residential-area column headers (row 10) ----
$ws.Range('R10').Value = 'CVA'
$ws.Range('S10').Value = 'RA'
$ws.Range('T10').Value = 'PA'
$ws.Range('U10').Value = 'SA'
$ws.Range('V10').Value = 'EA'

# ---- Miles distances from existing nodes (rows 11-18) to new residential nodes (cols R-V) ----
$ws.Range('R11').Value = 2.1
$ws.Range('S11').Value = 1.2
$ws.Range('T11').Value = 1
$ws.Range('U11').Value = 3.1
$ws.Range('V11').Value = 3.1
$ws.Range('R12').Value = 2.8
$ws.Range('S12').Value = 3.5
$ws.Range('T12').Value = 3.7
$ws.Range('U12').Value = 2.8
$ws.Range('V12').Value = 1.1
$ws.Range('R13').Value = 4.5
$ws.Range('S13').Value = 6.3
$ws.Range('T13').Value = 6.5
$ws.Range('U13').Value = 3.8
$ws.Range('V13').Value = 4.3
$ws.Range('R14').Value = 1.9
$ws.Range('S14').Value = 0.6
$ws.Range('T14').Value = 0.4
$ws.Range('U14').Value = 3
$ws.Range('V14').Value = 3
$ws.Range('R15').Value = 4.1
$ws.Range('S15').Value = 3.5
$ws.Range('T15').Value = 3.7
$ws.Range('U15').Value = 2.8
$ws.Range('V15').Value = 1.1
$ws.Range('R16').Value = 3.7
$ws.Range('S16').Value = 3.2
$ws.Range('T16').Value = 3.4
$ws.Range('U16').Value = 2.8
$ws.Range('V16').Value = 1.1
$ws.Range('R17').Value = 1.9
$ws.Range('S17').Value = 2.5
$ws.Range('T17').Value = 2.3
$ws.Range('U17').Value = 2.4
$ws.Range('V17').Value = 2.6
$ws.Range('R18').Value = 4.6
$ws.Range('S18').Value = 6.3
$ws.Range('T18').Value = 6.5
$ws.Range('U18').Value = 3.7
$ws.Range('V18').Value = 4.3

# ---- New residential-area rows (19-23): row labels + full distance rows ----
$ws.Range('I19').Value = 'CVA'
$ws.Range('I20').Value = 'RA'
$ws.Range('I21').Value = 'PA'
$ws.Range('I22').Value = 'SA'
$ws.Range('I23').Value = 'EA'

$ws.Range('J19').Value = 2.1
$ws.Range('K19').Value = 2.8
$ws.Range('L19').Value = 4.5
$ws.Range('M19').Value = 1.9
$ws.Range('N19').Value = 4.1
$ws.Range('O19').Value = 3.7
$ws.Range('P19').Value = 1.9
$ws.Range('Q19').Value = 4.6
$ws.Range('R19').Value = 0
$ws.Range('S19').Value = 2.5
$ws.Range('T19').Value = 2.3
$ws.Range('U19').Value = 2.2
$ws.Range('V19').Value = 2.2
$ws.Range('J20').Value = 1.2
$ws.Range('K20').Value = 3.5
$ws.Range('L20').Value = 6.3
$ws.Range('M20').Value = 0.6
$ws.Range('N20').Value = 3.5
$ws.Range('O20').Value = 3.2
$ws.Range('P20').Value = 2.5
$ws.Range('Q20').Value = 6.3
$ws.Range('R20').Value = 2.5
$ws.Range('S20').Value = 0
$ws.Range('T20').Value = 0.2
$ws.Range('U20').Value = 4
$ws.Range('V20').Value = 3.5
$ws.Range('J21').Value = 1
$ws.Range('K21').Value = 3.7
$ws.Range('L21').Value = 6.5
$ws.Range('M21').Value = 0.4
$ws.Range('N21').Value = 3.7
$ws.Range('O21').Value = 3.4
$ws.Range('P21').Value = 2.3
$ws.Range('Q21').Value = 6.5
$ws.Range('R21').Value = 2.3
$ws.Range('S21').Value = 0.2
$ws.Range('T21').Value = 0
$ws.Range('U21').Value = 4.1
$ws.Range('V21').Value = 3.3
$ws.Range('J22').Value = 3.1
$ws.Range('K22').Value = 2.8
$ws.Range('L22').Value = 3.8
$ws.Range('M22').Value = 3
$ws.Range('N22').Value = 2.8
$ws.Range('O22').Value = 2.8
$ws.Range('P22').Value = 2.4
$ws.Range('Q22').Value = 3.7
$ws.Range('R22').Value = 2.2
$ws.Range('S22').Value = 4
$ws.Range('T22').Value = 4.1
$ws.Range('U22').Value = 0
$ws.Range('V22').Value = 2
$ws.Range('J23').Value = 3.1
$ws.Range('K23').Value = 1.1
$ws.Range('L23').Value = 4.3
$ws.Range('M23').Value = 3
$ws.Range('N23').Value = 1.1
$ws.Range('O23').Value = 1.1
$ws.Range('P23').Value = 2.6
$ws.Range('Q23').Value = 4.3
$ws.Range('R23').Value = 2.2
$ws.Range('S23').Value = 3.5
$ws.Range('T23').Value = 3.3
$ws.Range('U23').Value = 2
$ws.Range('V23').Value = 0


# ---- Selection / view bookkeeping to match the saved workbook state ----
$ws.Range('V27').Select()
